# Update countries & provincias Spain
# - Refresh the COVID-19 data pull: Bolivia overtook Panama, Honduras
#   overtook Guatemala (both now sorted ahead of their former neighbour),
#   plus a handful of other countries received refreshed totals.
# - Update the "last updated" timestamp shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43/44: Bolivia now ranks above Panama (was the reverse).
$ws.Range("A43").Value = "Bolivia"
$ws.Range("B43").Value = 38071
$ws.Range("C43").Value = 1253
$ws.Range("D43").Value = 11272
$ws.Range("E43").Value = 25421
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 58
$ws.Range("H43").Value = 1378

$ws.Range("A44").Value = "Panama"
$ws.Range("B44").Value = 36983
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 17761
$ws.Range("E44").Value = 18502
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 720

# Row 55/56: Honduras now ranks above Guatemala (was the reverse).
$ws.Range("A55").Value = "Honduras"
$ws.Range("B55").Value = 22921
$ws.Range("C55").Value = 805
$ws.Range("D55").Value = 2387
$ws.Range("E55").Value = 19905
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 24
$ws.Range("H55").Value = 629

$ws.Range("A56").Value = "Guatemala"
$ws.Range("B56").Value = 22501
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 3330
$ws.Range("E56").Value = 18251
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 920

# Row 75: Australia refreshed totals (no reordering).
$ws.Range("B75").Value = 8443
$ws.Range("C75").Value = 81
$ws.Range("D75").Value = 7399
$ws.Range("E75").Value = 940

# Row 83: Venezuela refreshed totals (no reordering).
$ws.Range("B83").Value = 6750
$ws.Range("E83").Value = 4588
$ws.Range("H83").Value = 62

# Row 125: Nueva Zelanda refreshed totals (no reordering).
$ws.Range("B125").Value = 1533
$ws.Range("C125").Value = 3
$ws.Range("E125").Value = 21

# Fiyi now listed ahead of Dominica in the country list (tied case counts,
# so only the shared-string ordering changes; the data rows stay put).
$ws.Range("A205").Value = "Fiyi"
$ws.Range("A206").Value = "Dominica"

# Refresh "last updated" timestamp.
$ws.Range("A1").Value = "Datos actualizados a 5 de Julio de 2020 a las 05:11"
